# Atualização de bases das ligas, do dia: 14-06-2024 às 20:31
#
# This script swaps the match-record data between specific row pairs in the
# "Costa Rica Primera Division" sheet. Each pair corresponds to two fixtures
# that were entered in the wrong order (wrong id/row); the fix re-associates
# every data column (B through AD) between the two rows while leaving the
# row's own sequential index in column A untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data columns that can carry match information (everything except the
# running index in column A).
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

# Row pairs whose data needs to be swapped.
$rowPairs = @(
    @(95, 96),
    @(231, 232),
    @(237, 238),
    @(267, 269)
)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    foreach ($col in $cols) {
        $cellA = $ws.Range($col + $rowA)
        $cellB = $ws.Range($col + $rowB)

        $valA = $cellA.Value()
        $valB = $cellB.Value()

        if ($valA -ne $valB) {
            $cellA.Value = $valB
            $cellB.Value = $valA
        }
    }
}
